# Posthuma_2018_Table1.xlsx — "Change p to e-10 notation"
#
# On the "Table1_melted_phase3" sheet, the P-value column (F) was
# rewritten from rich-text "x.xx × 10^-y" superscript notation into
# plain "x.xxEy" scientific notation (left-aligned, non-bold/underline,
# scientific number format). One stray P value living in column J
# (row 14) got the same notation swap without any style changes. Two
# P-values that were already stored as plain numbers (F8, F27) only
# pick up the new number format/alignment. A new column G was opened
# up (widths on F/G set to 22), a spacer row 35 appended, row heights
# bumped from 18 to 19 on the affected rows, and the selection moved
# to G8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table1_melted_phase3")
$ws.Activate()

# Unicode glyphs used by the source data (minus sign vs. en dash —
# the original authors were inconsistent, so we reproduce it exactly).
$MINUS  = [char]0x2212
$ENDASH = [char]0x2013

# --- 1. Header cell F1: keep its font, just left-align it -----------------
$f1 = $ws.Range("F1")
$f1.HorizontalAlignment = -4131   # xlLeft

# --- 2. Data cells: rewrite P-value text, restyle to the new look ---------
$edits = @(
    @{Cell="F2"; Value="2.05E" + $MINUS + "10"},
    @{Cell="F3"; Value="1.10E" + $MINUS + "18"},
    @{Cell="F4"; Value="3.38E" + $MINUS + "44"},
    @{Cell="F5"; Value="8.92E" + $MINUS + "10"},
    @{Cell="F6"; Value="1.24E" + $MINUS + "8"},
    @{Cell="F7"; Value="1.93E" + $MINUS + "9"},
    @{Cell="F9"; Value="8.41E" + $MINUS + "11"},
    @{Cell="F10"; Value="1.45E" + $MINUS + "16"},
    @{Cell="F11"; Value="2.52E" + $MINUS + "10"},
    @{Cell="F12"; Value="2.22E" + $MINUS + "15"},
    @{Cell="F13"; Value="3.59E" + $MINUS + "11"},
    @{Cell="F14"; Value="2.10E" + $MINUS + "9"},
    @{Cell="F15"; Value="2.61E" + $MINUS + "19"},
    @{Cell="F16"; Value="1.26E" + $MINUS + "8"},
    @{Cell="F17"; Value="1.55E" + $MINUS + "15"},
    @{Cell="F18"; Value="2.19E" + $MINUS + "18"},
    @{Cell="F19"; Value="1.09E" + $MINUS + "11"},
    @{Cell="F20"; Value="1.65E" + $MINUS + "10"},
    @{Cell="F21"; Value="1.31E" + $MINUS + "9"},
    @{Cell="F22"; Value="3.35E" + $MINUS + "8"},
    @{Cell="F23"; Value="3.98E" + $MINUS + "8"},
    @{Cell="F24"; Value="9.16E" + $MINUS + "10"},
    @{Cell="F25"; Value="1.87E" + $MINUS + "8"},
    @{Cell="F26"; Value="9.66E" + $MINUS + "7"},
    @{Cell="F28"; Value="3.30E" + $MINUS + "8"},
    @{Cell="F29"; Value="7.93E" + $MINUS + "11"},
    @{Cell="F30"; Value="5.79E" + $MINUS + "276"},
    @{Cell="F31"; Value="4.64E" + $MINUS + "8"},
    @{Cell="F32"; Value="6.34E" + $MINUS + "9"},
    @{Cell="F33"; Value="6.56E" + $MINUS + "10"}
)

foreach ($edit in $edits) {
    $cell = $ws.Range($edit.Cell)
    $cell.NumberFormat = "0.00E+00"
    $cell.HorizontalAlignment = -4131   # xlLeft
    $cell.Font.Bold = $false
    $cell.Font.Underline = $false
    $cell.Value = $edit.Value
}

# F8 / F27 were already plain numbers (not rich text) — only the display
# format/alignment/font change, the underlying numeric value is untouched.
foreach ($addr in @("F8", "F27")) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "0.00E+00"
    $cell.HorizontalAlignment = -4131   # xlLeft
    $cell.Font.Bold = $false
    $cell.Font.Underline = $false
}

# J14 carries its own stray P-value; only the text changes, its style
# (plain, unaligned) is left alone.
$ws.Range("J14").Value = "2.59E" + $ENDASH + "4"

# --- 3. Row heights: 18 -> 19 on every affected data row -------------------
$rowsToResize = @(2,3,4,5,6,7,9,10,11,12,13,15,16,17,18,19,20,21,22,23,24,25,28,29,30,31,32,33)
foreach ($r in $rowsToResize) {
    $ws.Rows.Item($r).RowHeight = 19
}

# --- 4. New columns F & G sized to 22 characters ---------------------------
$ws.Columns.Item(6).ColumnWidth = 21 + 1/6   # renders as width="22" in xlsx
$ws.Columns.Item(7).ColumnWidth = 21 + 1/6

# --- 5. Spacer row 35 with a blank, pre-formatted F cell -------------------
$f35 = $ws.Range("F35")
$f35.NumberFormat = "0.00E+00"
$f35.HorizontalAlignment = -4131   # xlLeft

# --- 6. Selection moves to G8 ----------------------------------------------
$ws.Range("G8").Select()
